$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B3 ("Responsável Técnico" for the API/Tagging Coverage row) is
# renamed from "Yuri" to "Alexandre do Carmo".
$ws.Range("B3").Value = "Alexandre do Carmo"

# The saved view's active cell moved from K11 to B4.
[void]$ws.Range("B4").Select()
